$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: RowNum, Fecha(D), Variedad(H), Calidad(I), Volumen(J), PrecioMin(K), PrecioMax(L), PrecioProm(M), Origen(O), PrecioKg(P)
$rows = @(
    @(428, 44939, 'Crespo record', 'Primera', 2000, 800, 900, 850, 'Región Metropolitana', 850),
    @(429, 44939, 'Crespo record', 'Segunda', 1000, 700, 700, 700, 'Región Metropolitana', 700),
    @(430, 44939, 'Morada(o)', 'Primera', 1000, 1100, 1200, 1150, 'Región Metropolitana', 1150),
    @(431, 44939, 'Morada(o)', 'Segunda', 500, 800, 800, 800, 'Región Metropolitana', 800),
    @(432, 44425, 'Crespo record', 'Primera', 2000, 600, 700, 650, 'Región Metropolitana', 650),
    @(433, 44425, 'Crespo record', 'Segunda', 1000, 500, 500, 500, 'Región Metropolitana', 500),
    @(434, 44467, 'Crespo record', 'Primera', 1000, 600, 700, 650, 'Región Metropolitana', 650),
    @(435, 44467, 'Crespo record', 'Segunda', 500, 500, 500, 500, 'Región Metropolitana', 500),
    @(436, 44707, 'Crespo record', 'Primera', 2000, 1000, 1100, 1050, 'Región Metropolitana', 1050),
    @(437, 44707, 'Crespo record', 'Segunda', 1000, 800, 800, 800, 'Región Metropolitana', 800),
    @(438, 44707, 'Morada(o)', 'Primera', 1000, 1500, 1600, 1550, 'Región Metropolitana', 1550),
    @(439, 44488, 'Crespo record', 'Primera', 1000, 800, 900, 850, 'Región Metropolitana', 850),
    @(440, 44488, 'Crespo record', 'Segunda', 500, 700, 700, 700, 'Región Metropolitana', 700),
    @(441, 44504, 'Copenhague', 'Primera', 1800, 700, 850, 783, 'Región del Maule', 783),
    @(442, 44370, 'Crespo record', 'Primera', 800, 600, 700, 650, 'Región Metropolitana', 650),
    @(443, 44370, 'Crespo record', 'Segunda', 400, 500, 500, 500, 'Región Metropolitana', 500),
    @(444, 44691, 'Crespo record', 'Primera', 2000, 1100, 1200, 1150, 'Región Metropolitana', 1150),
    @(445, 44691, 'Crespo record', 'Segunda', 1000, 900, 900, 900, 'Región Metropolitana', 900),
    @(446, 44306, 'Crespo record', 'Primera', 1000, 700, 800, 750, 'Región Metropolitana', 750),
    @(447, 44306, 'Crespo record', 'Segunda', 500, 600, 600, 600, 'Región Metropolitana', 600),
    @(448, 44356, 'Crespo record', 'Primera', 1000, 600, 700, 650, 'Región Metropolitana', 650),
    @(449, 44356, 'Crespo record', 'Segunda', 500, 500, 500, 500, 'Región Metropolitana', 500),
    @(450, 44812, 'Crespo record', 'Primera', 1000, 1300, 1400, 1350, 'Región Metropolitana', 1350),
    @(451, 44812, 'Crespo record', 'Segunda', 500, 900, 900, 900, 'Región Metropolitana', 900),
    @(452, 44791, 'Copenhague', 'Primera', 1000, 1400, 1500, 1450, 'Región Metropolitana', 1450),
    @(453, 44335, 'Copenhague', 'Primera', 800, 700, 800, 750, 'Región Metropolitana', 750),
    @(454, 44335, 'Copenhague', 'Segunda', 400, 600, 600, 600, 'Región Metropolitana', 600),
    @(455, 44335, 'Crespo record', 'Primera', 800, 600, 700, 650, 'Región Metropolitana', 650),
    @(456, 44335, 'Crespo record', 'Segunda', 400, 500, 500, 500, 'Región Metropolitana', 500),
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = 11
    $ws.Cells.Item($rowNum, 2).Value = 'Vega Monumental Concepción'
    $ws.Cells.Item($rowNum, 3).Value = 'Bíobío'
    $ws.Cells.Item($rowNum, 4).Value = $r[1]
    $ws.Cells.Item($rowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($rowNum, 5).Value = 8
    $ws.Cells.Item($rowNum, 6).Value = 100112006
    $ws.Cells.Item($rowNum, 7).Value = 'Repollo'
    $ws.Cells.Item($rowNum, 8).Value = $r[2]
    $ws.Cells.Item($rowNum, 9).Value = $r[3]
    $ws.Cells.Item($rowNum, 10).Value = $r[4]
    $ws.Cells.Item($rowNum, 11).Value = $r[5]
    $ws.Cells.Item($rowNum, 12).Value = $r[6]
    $ws.Cells.Item($rowNum, 13).Value = $r[7]
    $ws.Cells.Item($rowNum, 14).Value = '$/unidad'
    $ws.Cells.Item($rowNum, 15).Value = $r[8]
    $ws.Cells.Item($rowNum, 16).Value = $r[9]
    $ws.Cells.Item($rowNum, 17).Value = 1
    $ws.Cells.Item($rowNum, 18).Value = 'Hortaliza'
}
